# Update "想去人数" (want-to-go count) values in column F across sheets.
# This mirrors a data refresh commit ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 254
$ws1.Range("F7").Value  = 616
$ws1.Range("F10").Value = 363
$ws1.Range("F11").Value = 159
$ws1.Range("F12").Value = 720
$ws1.Range("F13").Value = 97
$ws1.Range("F14").Value = 1841
$ws1.Range("F15").Value = 373
$ws1.Range("F16").Value = 3994
$ws1.Range("F17").Value = 378
$ws1.Range("F19").Value = 10
$ws1.Range("F20").Value = 62

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 482
$ws2.Range("F14").Value = 41

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5373
$ws3.Range("F3").Value = 333
$ws3.Range("F4").Value = 293

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5373
$ws4.Range("F4").Value  = 333
$ws4.Range("F6").Value  = 293
$ws4.Range("F7").Value  = 254
$ws4.Range("F12").Value = 482
$ws4.Range("F18").Value = 616
$ws4.Range("F22").Value = 363
$ws4.Range("F23").Value = 159
$ws4.Range("F26").Value = 720
$ws4.Range("F27").Value = 97
$ws4.Range("F29").Value = 1841
$ws4.Range("F30").Value = 373
$ws4.Range("F31").Value = 3994
$ws4.Range("F32").Value = 41
$ws4.Range("F33").Value = 378
$ws4.Range("F35").Value = 10
$ws4.Range("F36").Value = 62
